$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.073.21"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.71%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.741.50"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.38%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.07%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'404.90"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -4.22%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'128.38"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -2.43%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'3.726.38"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.24%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.604"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -5.53%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'1.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -0.09%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.717"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -5.46%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -7.91%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.0000358"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -4.01%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'40.37"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -5.03%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'4.323.69"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +0.38%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'9.64"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -5.77%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'14.41"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +11.13%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  -1.79%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.751.50"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +1.14%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  -7.26%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'66.206.96"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.70%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  -5.96%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'407.73"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -8.32%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'14.47"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -7.08%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'85.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "'3.02"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -4.92%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'36.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -4.97%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'5.58"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +11.87%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  -6.49%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'9.34"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -7.99%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'12.37"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -1.29%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -1.98%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -3.48%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'7.08"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -2.40%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'0.155"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -5.18%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'38.60"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -8.27%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  +0.00%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'55.13"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -2.32%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.0₃0724"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -0.76%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  -7.03%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'2.84"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -6.51%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.996"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -0.19%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  -7.64%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'3.16"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +20.85%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'145.18"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.00%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'26.41"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -7.72%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'2.05"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -5.80%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'3.22"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -7.70%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("B48").Value = "'Stacks"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'2.82"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -2.80%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("B49").Value = "'NEARProtocol"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = "'4.21"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -3.49%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'2.53"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -4.71%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.288"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -5.69%  "
$ws.Range("E51").Style = "Normal"
